$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Placas")

# Row 2: plate text corrected, schedule ("Escala") changed to "Plantonista 7 - 19"
$ws.Range("A2").Value = "RUN5C52"
$ws.Range("B2").Value = "Plantonista 7 - 19"

# Row 3: new plate added with same schedule
$ws.Range("A3").Value = "RUN5B64"
$ws.Range("B3").Value = "Plantonista 7 - 19"
